$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.319.58"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.919.82"
$ws.Range("E3").Value = "  -0.32%  "
$ws.Range("E4").Value = "  +0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "486.06"
$ws.Range("E5").Value = "  +1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.09"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  +0.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.997"
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.738"
$ws.Range("E9").Value = "  +2.11%  "
$ws.Range("E10").Value = "  +2.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000349"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.27"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.78"
$ws.Range("E13").Value = "  +4.73%  "
$ws.Range("D14").Value = "4.545.85"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "3.917.84"
$ws.Range("E15").Value = "  -1.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.29"
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.22"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D20").Value = "68.319.43"
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "432.02"
$ws.Range("E21").Value = "  -0.47%  "
$ws.Range("E22").Value = "  +7.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.17"
$ws.Range("E23").Value = "  +6.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.70"
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.65"
$ws.Range("E25").Value = "  +15.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.73"
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.25"
$ws.Range("E27").Value = "  +13.14%  "
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "719.19"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.79"
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.131"
$ws.Range("E32").Value = "  +3.12%  "
$ws.Range("E33").Value = "  +4.73%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0896"
$ws.Range("E34").Value = "  +4.96%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("E35").Value = "  +16.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.01"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "60.81"
$ws.Range("E37").Value = "  -4.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.404"
$ws.Range("E38").Value = "  +21.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.08"
$ws.Range("E39").Value = "  +11.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0497"
$ws.Range("E41").Value = "  +6.75%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.144"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.14"
$ws.Range("E43").Value = "  +4.11%  "
$ws.Range("E44").Value = "  +2.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.143"
$ws.Range("E45").Value = "  +1.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.34"
$ws.Range("E46").Value = "  +5.29%  "
$ws.Range("E47").Value = "  +0.44%  "
$ws.Range("E48").Value = "  +1.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.15"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.26"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").Value = "0.0₆0334"
$ws.Range("E51").Value = "  +32.18%  "
